# Auto-generated edit script applying scheduled market-data refresh
# to the Maduin_Profits workbook (per-sheet crafting-leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 436.21875
$ws.Range("I15").Value = 436.21875
$ws.Range("M15").Value = -1139.65625
$ws.Range("K15").Value = 1308.65625
$ws.Range("J19").Value = 421
$ws.Range("H19").Value = 491.7143
$ws.Range("L19").Value = 421
$ws.Range("N19").Value = -771
$ws.Range("I33").Value = 176
$ws.Range("H33").Value = 135.11111
$ws.Range("M33").Value = 53
$ws.Range("K33").Value = 176
$ws.Range("I41").Value = 765
$ws.Range("N41").Value = -1616
$ws.Range("M41").Value = -325
$ws.Range("H41").Value = 754.4545000000001
$ws.Range("J41").Value = 736
$ws.Range("L41").Value = 736
$ws.Range("K41").Value = 765
$ws.Range("N88").Value = -2799.3334
$ws.Range("L88").Value = 1987.3334
$ws.Range("J88").Value = 1987.3334
$ws.Range("H88").Value = 1683
$ws.Range("J91").Value = 1987.3334
$ws.Range("L91").Value = 1987.3334
$ws.Range("N91").Value = -4795.3334
$ws.Range("H91").Value = 1683
$ws.Range("K100").Value = 2692.9
$ws.Range("L100").Value = 0
$ws.Range("I100").Value = 2692.9
$ws.Range("J100").Value = 0
$ws.Range("M100").Value = -2151.9
$ws.Range("H100").Value = 2692.9
$ws.Range("M111").Value = -1681.25
$ws.Range("J111").Value = 2649.5
$ws.Range("N111").Value = -14082.5
$ws.Range("I111").Value = 1582.75
$ws.Range("K111").Value = 4748.25
$ws.Range("L111").Value = 7948.5
$ws.Range("H111").Value = 1796.1
$ws.Range("I125").Value = 0
$ws.Range("H125").Value = 6445
$ws.Range("K125").Value = 0
$ws.Range("N125").Value = -62925
$ws.Range("J125").Value = 6445
$ws.Range("L125").Value = 58005
$ws.Range("K131").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("L132").Value = 10888.5
$ws.Range("H132").Value = 3577.6667
$ws.Range("N132").Value = -15948.5
$ws.Range("J132").Value = 3629.5
$ws.Range("N100").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 40000
$ws.Range("H32").Value = 8472.412
$ws.Range("L32").Value = 40000
$ws.Range("N32").Value = -40574
$ws.Range("H104").Value = 100013500
$ws.Range("J104").Value = 100013500
$ws.Range("L104").Value = 100013500
$ws.Range("N104").Value = -100020488
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("I132").Value = 1101.6
$ws.Range("H132").Value = 1120.2188
$ws.Range("K132").Value = 3304.8
$ws.Range("M132").Value = -774.7999999999997
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 69999
$ws.Range("J61").Value = 69999
$ws.Range("L61").Value = 69999
$ws.Range("N61").Value = -70625
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988
$ws.Range("M105").Value = -1834.625
$ws.Range("J105").Value = 6250
$ws.Range("L105").Value = 6250
$ws.Range("I105").Value = 3581.625
$ws.Range("K105").Value = 3581.625
$ws.Range("H105").Value = 4115.3
$ws.Range("N105").Value = -9744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("N7").Value = -3091.75
$ws.Range("L7").Value = 2865.75
$ws.Range("M7").Value = -1024.1818
$ws.Range("K7").Value = 1137.1818
$ws.Range("J7").Value = 2865.75
$ws.Range("H7").Value = 1598.1333
$ws.Range("I7").Value = 1137.1818
$ws.Range("I58").Value = 2065.25
$ws.Range("J58").Value = 2983.2
$ws.Range("K58").Value = 2065.25
$ws.Range("L58").Value = 2983.2
$ws.Range("H58").Value = 2575.2222
$ws.Range("N58").Value = -3389.2
$ws.Range("M58").Value = -1862.25
$ws.Range("I80").Value = 0
$ws.Range("H80").Value = 25000
$ws.Range("K80").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("H83").Value = 25000
$ws.Range("H104").Value = 41788.332
$ws.Range("J104").Value = 60182.5
$ws.Range("L104").Value = 60182.5
$ws.Range("N104").Value = -65424.5
$ws.Range("M105").Value = 551
$ws.Range("J105").Value = 1497
$ws.Range("L105").Value = 1497
$ws.Range("I105").Value = 1196
$ws.Range("K105").Value = 1196
$ws.Range("H105").Value = 1396.6666
$ws.Range("N105").Value = -4991
$ws.Range("I122").Value = 2330.5
$ws.Range("M122").Value = -4541.5
$ws.Range("K122").Value = 6991.5
$ws.Range("H122").Value = 2557.3333
$ws.Range("L136").Value = 8949.599999999999
$ws.Range("J136").Value = 2983.2
$ws.Range("H136").Value = 2575.2222
$ws.Range("I136").Value = 2065.25
$ws.Range("K136").Value = 6195.75
$ws.Range("N136").Value = -14049.6
$ws.Range("M136").Value = -3645.75
$ws.Range("M4").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 285714880
$ws.Range("L4").Value = 2376
$ws.Range("J4").Value = 792
$ws.Range("M4").Value = -1000001588
$ws.Range("I4").Value = 333333900
$ws.Range("N4").Value = -2600
$ws.Range("K4").Value = 1000001700
$ws.Range("K86").Value = 1135.8
$ws.Range("H86").Value = 1147.1666
$ws.Range("M86").Value = 50.19999999999982
$ws.Range("J86").Value = 4990
$ws.Range("L86").Value = 14970
$ws.Range("N86").Value = -17342
$ws.Range("I86").Value = 378.6
$ws.Range("L89").Value = 44910
$ws.Range("I89").Value = 378.6
$ws.Range("J89").Value = 4990
$ws.Range("M89").Value = 2520.6
$ws.Range("N89").Value = -56766
$ws.Range("K89").Value = 3407.4
$ws.Range("H89").Value = 1147.1666
$ws.Range("I105").Value = 7000
$ws.Range("K105").Value = 21000
$ws.Range("H105").Value = 7000
$ws.Range("M105").Value = -18379
$ws.Range("I122").Value = 1044.4286
$ws.Range("M122").Value = -6949.857399999999
$ws.Range("J122").Value = 1394.6
$ws.Range("K122").Value = 9399.857399999999
$ws.Range("L122").Value = 12551.4
$ws.Range("N122").Value = -17451.4
$ws.Range("H122").Value = 1190.3334
$ws.Range("H128").Value = 331570.8
$ws.Range("K128").Value = 994712.3999999999
$ws.Range("M128").Value = -989732.3999999999
$ws.Range("I128").Value = 331570.8
$ws.Range("K131").Value = 982.5
$ws.Range("H131").Value = 725
$ws.Range("M131").Value = 4057.5
$ws.Range("I131").Value = 327.5
$ws.Range("L132").Value = 6750
$ws.Range("H132").Value = 700
$ws.Range("N132").Value = -11810
$ws.Range("J132").Value = 750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L3").Value = 1816
$ws.Range("J3").Value = 1816
$ws.Range("K3").Value = 3333500
$ws.Range("M3").Value = -3333384
$ws.Range("I3").Value = 3333500
$ws.Range("N3").Value = -2048
$ws.Range("H3").Value = 1429680.6
$ws.Range("I122").Value = 1298.7778
$ws.Range("M122").Value = -1446.3334
$ws.Range("K122").Value = 3896.3334
$ws.Range("H122").Value = 2335.3635
$ws.Range("I132").Value = 3010
$ws.Range("H132").Value = 3402.5264
$ws.Range("K132").Value = 9030
$ws.Range("M132").Value = -6500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("M22").Value = -1855
$ws.Range("I22").Value = 2150
$ws.Range("K22").Value = 2150
$ws.Range("H22").Value = 2120
$ws.Range("K27").Value = 2150
$ws.Range("H27").Value = 2120
$ws.Range("I27").Value = 2150
$ws.Range("M27").Value = -2043
$ws.Range("H40").Value = 6146.4
$ws.Range("M40").Value = -5705.1816
$ws.Range("K40").Value = 5841.1816
$ws.Range("I40").Value = 5841.1816
$ws.Range("J81").Value = 0
$ws.Range("H81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("L132").Value = 23330.0772
$ws.Range("H132").Value = 6095.269
$ws.Range("N132").Value = -28390.0772
$ws.Range("J132").Value = 7776.6924
$ws.Range("M2").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("H136").Value = 735.625
$ws.Range("I136").Value = 735.625
$ws.Range("K136").Value = 2206.875
$ws.Range("M136").Value = 343.125
$ws.Range("M2").ClearContents()
